# Update the "想去人数" (interested-attendee count) figures in column F
# for the "展览" (Exhibitions) and "全部类型" (All types) worksheets.
# These two sheets list overlapping events, but "全部类型" has one extra
# row (a 演出/show entry inserted at row 7), so the same events live one
# row lower there than in "展览".

$wb = $excel.ActiveWorkbook

# Sheet "展览": rows 2,3,6,8,10,11,12
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value = 94
$wsExpo.Range("F3").Value = 817
$wsExpo.Range("F6").Value = 129
$wsExpo.Range("F8").Value = 4788
$wsExpo.Range("F10").Value = 5143
$wsExpo.Range("F11").Value = 589
$wsExpo.Range("F12").Value = 1286

# Sheet "全部类型": same events, shifted down one row (2,3,6,9,11,12,13)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 94
$wsAll.Range("F3").Value = 817
$wsAll.Range("F6").Value = 129
$wsAll.Range("F9").Value = 4788
$wsAll.Range("F11").Value = 5143
$wsAll.Range("F12").Value = 589
$wsAll.Range("F13").Value = 1286
